# Apply the updates described in the diff:
#  - Metadata sheet: URL and Date values change
#  - Include ValueSets sheet: ValueSet URL changes
#  - Include ValueSets 2 sheet: ValueSet URL changes

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "http://example.org/ig/exampleig/ValueSet/presence-valueset"
$wsMeta.Range("B8").Value = "2023-04-27T11:07:01-05:00"

$wsSct = $wb.Worksheets.Item("Include ValueSets")
$wsSct.Range("A2").Value = "http://example.org/ig/exampleig/ValueSet/presence-sct-valueset"

$wsLnc = $wb.Worksheets.Item("Include ValueSets 2")
$wsLnc.Range("A2").Value = "http://example.org/ig/exampleig/ValueSet/presence-lnc-valueset"
